# Popup handling at profile management
#
# The "Manage Products" sheet has a popup used to (re)submit a product's
# name. Submitting it for the first six products (rows 2-7, column B)
# replaces the placeholder product-name text that was sitting in each
# cell with a freshly generated one - the cell's existing formatting
# (fill/border) is left exactly as-is, only the text changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "prodXkLk"
$ws.Range("B3").Value = "proddASw"
$ws.Range("B4").Value = "prodACZF"
$ws.Range("B5").Value = "prodcMOg"
$ws.Range("B6").Value = "prodAmgZ"
$ws.Range("B7").Value = "prodUjkK"
